$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# The sheet has 5 existing 6-column "period" blocks (B:G, H:M, N:S, T:Y,
# Z:AE), each showing, for one period (Mar 2017 .. Mar 2021):
#   col1: Shareholder Funds (number)      col2: Debts (text)
#   col3: Total Revenue (number)          col4: PBT (text)
#   col5: PAT (text)                      col6: Cash & Cash Eq (header only)
# This change appends a 6th block (AF:AK) for "Mar 2022", mirroring the
# layout/formatting of the preceding ("Z:AE" / Mar 2021) block, with new
# data: Shareholder Funds 26.15, Debts 1.25, Total Revenue 31.89,
# PBT 3.77, PAT 2.67.
# ----------------------------------------------------------------------

# 1) Clone the formatting (font/alignment/style) of the most recent block
#    (Z:AE, rows 1-4) onto the new block (AF:AK).
$ws.Range("Z1:AE4").Copy()
$ws.Range("AF1:AK4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 2) Merge the new block's header cell, same as the other period headers.
$ws.Range("AF1:AK1").Merge()

# 3) Header label for the new block. Force text format first so Excel
#    doesn't auto-convert the "Mar 2022" text into a date serial, then
#    re-apply the clean header look from the neighboring header cell.
$ws.Range("AF1").NumberFormat = "@"
$ws.Range("AF1").Value2 = "Mar 2022"
$ws.Range("Z1").Copy()
$ws.Range("AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 4) Row 2/3 labels (plain text - no coercion risk).
$ws.Range("AF2").Value2 = "Shareholder"
$ws.Range("AG2").Value2 = "Debts"
$ws.Range("AH2").Value2 = "Total"
$ws.Range("AI2").Value2 = "PBT"
$ws.Range("AJ2").Value2 = "PAT"
$ws.Range("AK2").Value2 = "Cash"

$ws.Range("AF3").Value2 = "Funds"
$ws.Range("AH3").Value2 = "Revenue"
$ws.Range("AK3").Value2 = "Cash Eq"

# 5) Row 4 data for the new "Mar 2022" period.
#    Shareholder Funds & Total Revenue are numeric; Debts/PBT/PAT are
#    stored as text (matching every other period block), so again force
#    text format before assignment to dodge numeric auto-conversion, then
#    restore the clean data-row look from the analogous cell.
$ws.Range("AF4").Value2 = 26.15

$ws.Range("AG4").NumberFormat = "@"
$ws.Range("AG4").Value2 = "1.25"
$ws.Range("AA4").Copy()
$ws.Range("AG4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("AH4").Value2 = 31.89

$ws.Range("AI4").NumberFormat = "@"
$ws.Range("AI4").Value2 = "3.77"
$ws.Range("AC4").Copy()
$ws.Range("AI4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("AJ4").NumberFormat = "@"
$ws.Range("AJ4").Value2 = "2.67"
$ws.Range("AD4").Copy()
$ws.Range("AJ4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Output "Mar 2022 block added"
